$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Format B1: bold font, thin border on all 4 sides, centered horizontally, top vertically
$c1 = $ws.Range("B1")
$c1.Font.Bold = $true
$c1.HorizontalAlignment = -4108   # xlCenter
$c1.VerticalAlignment = -4160     # xlTop
$c1.Borders.LineStyle = 1         # xlContinuous
$c1.Borders.Weight = 2            # xlThin

# Apply the exact same formatting to A2 by copying B1's format over --
# this reuses the same style record instead of minting a near-duplicate one.
$c1.Copy()
$c2 = $ws.Range("A2")
$c2.PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = 0
